$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The activity rows were updated/shuffled as new match data rolled in.
# Ensure the cells keep their original "text" storage (values look numeric
# but are stored as text in the source workbook) by forcing a text format
# before writing the values.
$cells = @("C2","D2","E2","F2","C4","D4","E4","F4","C5","D5","E5","F5","C6","D6","E6","F6","C7","D7","E7","F7","C8","D8","E8","F8")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("C2").Value = "0"
$ws.Range("D2").Value = "3"
$ws.Range("E2").Value = "0"
$ws.Range("F2").Value = "0"

$ws.Range("C4").Value = "26"
$ws.Range("D4").Value = "25"
$ws.Range("E4").Value = "3"
$ws.Range("F4").Value = "0"

$ws.Range("C5").Value = "4"
$ws.Range("D5").Value = "3"
$ws.Range("E5").Value = "1"
$ws.Range("F5").Value = "0"

$ws.Range("C6").Value = "85"
$ws.Range("D6").Value = "42"
$ws.Range("E6").Value = "4"
$ws.Range("F6").Value = "7"

$ws.Range("C7").Value = "8"
$ws.Range("D7").Value = "9"
$ws.Range("E7").Value = "1"
$ws.Range("F7").Value = "0"

$ws.Range("C8").Value = "5"
$ws.Range("D8").Value = "9"
$ws.Range("E8").Value = "0"
$ws.Range("F8").Value = "0"
